$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the columns F:N (old data size no longer needed) across rows 1-5
$ws.Range("F1:N5").Clear()

# Row 1 - header values (only B1:E1 remain)
$ws.Range("B1").Value = 100
$ws.Range("C1").Value = 250
$ws.Range("D1").Value = 500
$ws.Range("E1").Value = 750

# Row 2 - Insertion Sort
$ws.Range("A2").Value = "Insertion Sort"
$ws.Range("B2").Value = 5.785
$ws.Range("C2").Value = 40.942
$ws.Range("D2").Value = 131.182
$ws.Range("E2").Value = 306.845

# Row 3 - Quicksort
$ws.Range("A3").Value = "Quicksort"
$ws.Range("B3").Value = 1.208
$ws.Range("C3").Value = 3.815
$ws.Range("D3").Value = 8.086
$ws.Range("E3").Value = 11.974

# Row 4 - Heap Sort
$ws.Range("A4").Value = "Heap Sort"
$ws.Range("B4").Value = 2.8
$ws.Range("C4").Value = 9.329000000000001
$ws.Range("D4").Value = 20.575
$ws.Range("E4").Value = 33.546

# Row 5 - Bucket Sort
$ws.Range("A5").Value = "Bucket Sort"
$ws.Range("B5").Value = 1.203
$ws.Range("C5").Value = 4.016
$ws.Range("D5").Value = 9.923
$ws.Range("E5").Value = 14.175

# Row 6 - IntroSort (new row) - copy A5's formatting (border/bold/alignment) first
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "IntroSort"
$ws.Range("B6").Value = 1.018
$ws.Range("C6").Value = 4.445
$ws.Range("D6").Value = 9.755000000000001
$ws.Range("E6").Value = 15.748
